$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.75
$ws.Range("I2").Value = 2.7
$ws.Range("J2").Value = 3.6
$ws.Range("L2").Value = 3.6
$ws.Range("O2").Value = 1.5
$ws.Range("P2").Value = 2.5
$ws.Range("Q2").Value = 2.6
$ws.Range("R2").Value = 1.48
$ws.Range("W2").Value = 6.5
$ws.Range("X2").Value = 12
$ws.Range("Y2").Value = 11
$ws.Range("Z2").Value = 29
$ws.Range("AC2").Value = 6.5
$ws.Range("AF2").Value = 67
$ws.Range("AG2").Value = 6.5
$ws.Range("AH2").Value = 12
$ws.Range("AJ2").Value = 29
$ws.Range("AN2").Value = 4.5
$ws.Range("AO2").Value = 17
$ws.Range("AQ2").Value = 51
$ws.Range("AS2").Value = 301
$ws.Range("AW2").Value = 4.5
$ws.Range("AX2").Value = 17
$ws.Range("AY2").Value = 34
$ws.Range("BA2").Value = 101

# Row 3
$ws.Range("Q3").Value = 2.1
$ws.Range("R3").Value = 1.7

# Row 4
$ws.Range("G4").Value = 2.15
$ws.Range("I4").Value = 3.4
$ws.Range("J4").Value = 2.75
$ws.Range("L4").Value = 3.75
$ws.Range("Q4").Value = 1.93
$ws.Range("R4").Value = 1.93
$ws.Range("S4").Value = 1.4
$ws.Range("T4").Value = 2.75
$ws.Range("W4").Value = 8.5
$ws.Range("X4").Value = 11
$ws.Range("Y4").Value = 9
$ws.Range("AC4").Value = 10
$ws.Range("AE4").Value = 13
$ws.Range("AK4").Value = 26
$ws.Range("AL4").Value = 34
$ws.Range("AM4").Value = 201
$ws.Range("AS4").Value = 151
$ws.Range("AT4").Value = 2.75
$ws.Range("AW4").Value = 5.5
$ws.Range("AY4").Value = 26
$ws.Range("BA4").Value = 81

# Row 6
$ws.Range("G6").Value = 3.1
$ws.Range("H6").Value = 3
$ws.Range("I6").Value = 2.38
$ws.Range("K6").Value = 1.91
$ws.Range("L6").Value = 3.25
$ws.Range("M6").Value = 1.11
$ws.Range("N6").Value = 6.5
$ws.Range("O6").Value = 1.53
$ws.Range("P6").Value = 2.38
$ws.Range("Q6").Value = 2.7
$ws.Range("R6").Value = 1.44
$ws.Range("U6").Value = 2.2
$ws.Range("V6").Value = 1.62
$ws.Range("W6").Value = 7
$ws.Range("X6").Value = 13
$ws.Range("AE6").Value = 21
$ws.Range("AF6").Value = 81
$ws.Range("AI6").Value = 11
$ws.Range("AJ6").Value = 23
$ws.Range("AS6").Value = 351
$ws.Range("AV6").Value = 81
$ws.Range("BB6").Value = 301

# Row 7
$ws.Range("O7").Value = 1.36
$ws.Range("P7").Value = 3
$ws.Range("Q7").Value = 2.15
$ws.Range("R7").Value = 1.67

# Row 8
$ws.Range("G8").Value = 1.83
$ws.Range("H8").Value = 3
$ws.Range("I8").Value = 4.5
$ws.Range("J8").Value = 2.6
$ws.Range("L8").Value = 5
$ws.Range("M8").Value = 1.08
$ws.Range("N8").Value = 8
$ws.Range("W8").Value = 6
$ws.Range("AD8").Value = 6
$ws.Range("AN8").Value = 3.75
$ws.Range("AQ8").Value = 41
$ws.Range("AX8").Value = 26
$ws.Range("BA8").Value = 126
